$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# -----------------------------------------------------------------
# Cell values: header row + the 3 data rows (row 4 is brand new)
# -----------------------------------------------------------------
# Row 1
$ws.Range("A1").Value = "SL #"
$ws.Range("B1").Value = "MONTH"
$ws.Range("C1").Value = "Emp. NAME"
$ws.Range("D1").Value = "DOJ"
$ws.Range("E1").Value = "STATUS"
$ws.Range("F1").Value = "DESIGNATION"
$ws.Range("G1").Value = "DEPARTMENT"
$ws.Range("H1").Value = "GROSS"
$ws.Range("I1").Value = "Per Month"
$ws.Range("J1").Value = "Actual Per Month"
$ws.Range("K1").Value = "Actual Days"
$ws.Range("L1").Value = "Working Days"
$ws.Range("M1").Value = "BASIC"
$ws.Range("N1").Value = "DA"
$ws.Range("O1").Value = "HRA"
$ws.Range("P1").Value = "Spcl Allowance"
$ws.Range("Q1").Value = "Arrears"
$ws.Range("R1").Value = "Gross Pay"
$ws.Range("S1").Value = "PF"
$ws.Range("T1").Value = "ESIC"
$ws.Range("U1").Value = "PT"
$ws.Range("V1").Value = "TDS"
$ws.Range("W1").Value = "Deducted allowance1"
$ws.Range("X1").Value = "Deducted allowance2"
$ws.Range("Y1").Value = "total_deducations"
$ws.Range("Z1").Value = "NetPay"

# Row 2
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = "10-2014"
$ws.Range("C2").Value = "Sekhar Beri"
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "01/06/2014"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "New"
$ws.Range("F2").Value = "Senior HR"
$ws.Range("G2").Value = "HR"
$ws.Range("H2").Value = 120000.0
$ws.Range("I2").Value = 10000.0
$ws.Range("J2").Value = 9386.67
$ws.Range("K2").Value = 31.0
$ws.Range("L2").Value = 31.0
$ws.Range("M2").Value = 4000.0
$ws.Range("N2").Value = 800.0
$ws.Range("O2").Value = 1000.0
$ws.Range("P2").Value = 3461.67
$ws.Range("Q2").Value = 0.0
$ws.Range("R2").Value = 9386.67
$ws.Range("S2").Value = 480.0
$ws.Range("T2").Value = 164.27
$ws.Range("U2").Value = 0.0
$ws.Range("V2").Value = 0.0
$ws.Range("W2").Value = 0
$ws.Range("X2").Value = 0
$ws.Range("Y2").Value = 685.94
$ws.Range("Z2").Value = 8700.73

# Row 3
$ws.Range("A3").Value = 2
$ws.Range("B3").Value = "10-2014"
$ws.Range("C3").Value = "Priyanka Muddana"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "01/06/2014"
$ws.Range("D3").Style = "Normal"
$ws.Range("F3").Value = "Junior Development"
$ws.Range("G3").Value = "Development"
$ws.Range("H3").Value = 120000.0
$ws.Range("I3").Value = 10000.0
$ws.Range("J3").Value = 9561.67
$ws.Range("K3").Value = 31.0
$ws.Range("L3").Value = 31.0
$ws.Range("M3").Value = 4000.0
$ws.Range("N3").Value = 800.0
$ws.Range("O3").Value = 1000.0
$ws.Range("P3").Value = 3761.67
$ws.Range("Q3").Value = 0.0
$ws.Range("R3").Value = 9561.67
$ws.Range("S3").Value = 480.0
$ws.Range("T3").Value = 0.0
$ws.Range("U3").Value = 0.0
$ws.Range("V3").Value = 0.0
$ws.Range("W3").Value = 0
$ws.Range("X3").Value = 0
$ws.Range("Y3").Value = 521.67
$ws.Range("Z3").Value = 9040.0

# Row 4
$ws.Range("A4").Value = 3
$ws.Range("B4").Value = "10-2014"
$ws.Range("C4").Value = "Pattabhi RamaRao Galidevara"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "02/06/2014"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "Regular"
$ws.Range("F4").Value = "Junior Development"
$ws.Range("G4").Value = "Development"
$ws.Range("H4").Value = 120000.0
$ws.Range("I4").Value = 10000.0
$ws.Range("J4").Value = 10000.0
$ws.Range("K4").Value = 31.0
$ws.Range("L4").Value = 31.0
$ws.Range("M4").Value = 4000.0
$ws.Range("N4").Value = 0
$ws.Range("O4").Value = 0
$ws.Range("P4").Value = 6000.0
$ws.Range("Q4").Value = 0.0
$ws.Range("R4").Value = 10000.0
$ws.Range("S4").Value = 0.0
$ws.Range("T4").Value = 0.0
$ws.Range("U4").Value = 0.0
$ws.Range("V4").Value = 0.0
$ws.Range("W4").Value = 0
$ws.Range("X4").Value = 0
$ws.Range("Y4").Value = 0.0
$ws.Range("Z4").Value = 10000.0

# -----------------------------------------------------------------
# Blank cells (present in the source but emptied out)
# -----------------------------------------------------------------
$ws.Range("E3").ClearContents()

# -----------------------------------------------------------------
# Column widths - re-fit the columns whose best-fit width changed
# (values are the closest achievable ColumnWidth given the 1/7-char
# rounding Excel applies internally; matches target OOXML width)
# -----------------------------------------------------------------
$ws.Columns("C").ColumnWidth = 16.714285714285715
$ws.Columns("D").ColumnWidth = 13.428571428571429
$ws.Columns("I").ColumnWidth = 10.142857142857142
$ws.Columns("M").ColumnWidth = 7.857142857142857
$ws.Columns("O").ColumnWidth = 7.857142857142857
$ws.Columns("W").ColumnWidth = 15.571428571428571
$ws.Columns("X").ColumnWidth = 15.571428571428571
$ws.Columns("Y").ColumnWidth = 11.142857142857142
$ws.Columns("Z").ColumnWidth = 9.0

Write-Host "edit applied"
